$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Efna5"
$ws.Cells.Item(2,3).Value = "Ephb6"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 0.3227736666666667
$ws.Cells.Item(2,8).Value = 0.968321
$ws.Cells.Item(2,9).Value = 0.1416094457286952
$ws.Cells.Item(2,10).Value = 0.1416094457286952
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 0.3333333333333333
$ws.Cells.Item(2,13).Value = 0.04731666666666667
$ws.Cells.Item(2,14).Value = 0.14195
$ws.Cells.Item(2,15).Value = 0.01864859375546025
$ws.Cells.Item(2,16).Value = 0.01864859375546025
$ws.Cells.Item(2,17).Value = 0.01527257399444444
$ws.Cells.Item(2,18).Value = 0.13745316595
$ws.Cells.Item(2,19).Value = 0.002640817025330331
$ws.Cells.Item(2,20).Value = 0.002640817025330331

$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Efna5"
$ws.Cells.Item(3,3).Value = "Ephb6"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 0.3227736666666667
$ws.Cells.Item(3,8).Value = 0.968321
$ws.Cells.Item(3,9).Value = 0.1416094457286952
$ws.Cells.Item(3,10).Value = 0.1416094457286952
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 0.7341896666666666
$ws.Cells.Item(3,14).Value = 2.202569
$ws.Cells.Item(3,15).Value = 0.289361144764849
$ws.Cells.Item(3,16).Value = 0.289361144764849
$ws.Cells.Item(3,17).Value = 0.2369770907387778
$ws.Cells.Item(3,18).Value = 2.132793816649
$ws.Cells.Item(3,19).Value = 0.040976271325571
$ws.Cells.Item(3,20).Value = 0.040976271325571

$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Efna5"
$ws.Cells.Item(4,3).Value = "Ephb6"
$ws.Cells.Item(4,4).Value = "MuSCs"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 0.3227736666666667
$ws.Cells.Item(4,8).Value = 0.968321
$ws.Cells.Item(4,9).Value = 0.1416094457286952
$ws.Cells.Item(4,10).Value = 0.1416094457286952
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 1.683518666666667
$ws.Cells.Item(4,14).Value = 5.050556
$ws.Cells.Item(4,15).Value = 0.663513681459685
$ws.Cells.Item(4,16).Value = 0.663513681459685
$ws.Cells.Item(4,17).Value = 0.5433954929417778
$ws.Cells.Item(4,18).Value = 4.890559436476
$ws.Cells.Item(4,19).Value = 0.09395980466491199
$ws.Cells.Item(4,20).Value = 0.09395980466491199

$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Efna5"
$ws.Cells.Item(5,3).Value = "Ephb6"
$ws.Cells.Item(5,4).Value = "Resolving-Mac"
$ws.Cells.Item(5,5).Value = 2
$ws.Cells.Item(5,6).Value = 0.6666666666666666
$ws.Cells.Item(5,7).Value = 0.3227736666666667
$ws.Cells.Item(5,8).Value = 0.968321
$ws.Cells.Item(5,9).Value = 0.1416094457286952
$ws.Cells.Item(5,10).Value = 0.1416094457286952
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 0.072253
$ws.Cells.Item(5,14).Value = 0.216759
$ws.Cells.Item(5,15).Value = 0.02847658002000569
$ws.Cells.Item(5,16).Value = 0.02847658002000569
$ws.Cells.Item(5,17).Value = 0.02332136573766667
$ws.Cells.Item(5,18).Value = 0.209892291639
$ws.Cells.Item(5,19).Value = 0.00403255271288184
$ws.Cells.Item(5,20).Value = 0.004032552712881841

$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Efna5"
$ws.Cells.Item(6,3).Value = "Ephb6"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 1.874986333333333
$ws.Cells.Item(6,8).Value = 5.624959
$ws.Cells.Item(6,9).Value = 0.8226066833587575
$ws.Cells.Item(6,10).Value = 0.8226066833587576
$ws.Cells.Item(6,11).Value = 1
$ws.Cells.Item(6,12).Value = 0.3333333333333333
$ws.Cells.Item(6,13).Value = 0.04731666666666667
$ws.Cells.Item(6,14).Value = 0.14195
$ws.Cells.Item(6,15).Value = 0.01864859375546025
$ws.Cells.Item(6,16).Value = 0.01864859375546025
$ws.Cells.Item(6,17).Value = 0.0887181033388889
$ws.Cells.Item(6,18).Value = 0.79846293005
$ws.Cells.Item(6,19).Value = 0.01534045785848399
$ws.Cells.Item(6,20).Value = 0.01534045785848399

$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Efna5"
$ws.Cells.Item(7,3).Value = "Ephb6"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 1.874986333333333
$ws.Cells.Item(7,8).Value = 5.624959
$ws.Cells.Item(7,9).Value = 0.8226066833587575
$ws.Cells.Item(7,10).Value = 0.8226066833587576
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 0.7341896666666666
$ws.Cells.Item(7,14).Value = 2.202569
$ws.Cells.Item(7,15).Value = 0.289361144764849
$ws.Cells.Item(7,16).Value = 0.289361144764849
$ws.Cells.Item(7,17).Value = 1.376595591074556
$ws.Cells.Item(7,18).Value = 12.389360319671
$ws.Cells.Item(7,19).Value = 0.2380304115879058
$ws.Cells.Item(7,20).Value = 0.2380304115879058

$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Efna5"
$ws.Cells.Item(8,3).Value = "Ephb6"
$ws.Cells.Item(8,4).Value = "MuSCs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 1.874986333333333
$ws.Cells.Item(8,8).Value = 5.624959
$ws.Cells.Item(8,9).Value = 0.8226066833587575
$ws.Cells.Item(8,10).Value = 0.8226066833587576
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 1.683518666666667
$ws.Cells.Item(8,14).Value = 5.050556
$ws.Cells.Item(8,15).Value = 0.663513681459685
$ws.Cells.Item(8,16).Value = 0.663513681459685
$ws.Cells.Item(8,17).Value = 3.156574491911556
$ws.Cells.Item(8,18).Value = 28.409170427204
$ws.Cells.Item(8,19).Value = 0.5458107888687106
$ws.Cells.Item(8,20).Value = 0.5458107888687107

$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Efna5"
$ws.Cells.Item(9,3).Value = "Ephb6"
$ws.Cells.Item(9,4).Value = "Resolving-Mac"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 1.874986333333333
$ws.Cells.Item(9,8).Value = 5.624959
$ws.Cells.Item(9,9).Value = 0.8226066833587575
$ws.Cells.Item(9,10).Value = 0.8226066833587576
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 0.072253
$ws.Cells.Item(9,14).Value = 0.216759
$ws.Cells.Item(9,15).Value = 0.02847658002000569
$ws.Cells.Item(9,16).Value = 0.02847658002000569
$ws.Cells.Item(9,17).Value = 0.1354733875423333
$ws.Cells.Item(9,18).Value = 1.219260487881
$ws.Cells.Item(9,19).Value = 0.02342502504365714
$ws.Cells.Item(9,20).Value = 0.02342502504365715

$ws.Cells.Item(10,1).Value = "MuSCs"
$ws.Cells.Item(10,2).Value = "Efna5"
$ws.Cells.Item(10,3).Value = "Ephb6"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = 2
$ws.Cells.Item(10,6).Value = 0.6666666666666666
$ws.Cells.Item(10,7).Value = 0.081563
$ws.Cells.Item(10,8).Value = 0.244689
$ws.Cells.Item(10,9).Value = 0.03578387091254728
$ws.Cells.Item(10,10).Value = 0.03578387091254728
$ws.Cells.Item(10,11).Value = 1
$ws.Cells.Item(10,12).Value = 0.3333333333333333
$ws.Cells.Item(10,13).Value = 0.04731666666666667
$ws.Cells.Item(10,14).Value = 0.14195
$ws.Cells.Item(10,15).Value = 0.01864859375546025
$ws.Cells.Item(10,16).Value = 0.01864859375546025
$ws.Cells.Item(10,17).Value = 0.003859289283333333
$ws.Cells.Item(10,18).Value = 0.03473360355
$ws.Cells.Item(10,19).Value = 0.0006673188716459248
$ws.Cells.Item(10,20).Value = 0.0006673188716459248

$ws.Cells.Item(11,1).Value = "MuSCs"
$ws.Cells.Item(11,2).Value = "Efna5"
$ws.Cells.Item(11,3).Value = "Ephb6"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = 2
$ws.Cells.Item(11,6).Value = 0.6666666666666666
$ws.Cells.Item(11,7).Value = 0.081563
$ws.Cells.Item(11,8).Value = 0.244689
$ws.Cells.Item(11,9).Value = 0.03578387091254728
$ws.Cells.Item(11,10).Value = 0.03578387091254728
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 0.7341896666666666
$ws.Cells.Item(11,14).Value = 2.202569
$ws.Cells.Item(11,15).Value = 0.289361144764849
$ws.Cells.Item(11,16).Value = 0.289361144764849
$ws.Cells.Item(11,17).Value = 0.05988271178233333
$ws.Cells.Item(11,18).Value = 0.538944406041
$ws.Cells.Item(11,19).Value = 0.01035446185137226
$ws.Cells.Item(11,20).Value = 0.01035446185137226

$ws.Cells.Item(12,1).Value = "MuSCs"
$ws.Cells.Item(12,2).Value = "Efna5"
$ws.Cells.Item(12,3).Value = "Ephb6"
$ws.Cells.Item(12,4).Value = "MuSCs"
$ws.Cells.Item(12,5).Value = 2
$ws.Cells.Item(12,6).Value = 0.6666666666666666
$ws.Cells.Item(12,7).Value = 0.081563
$ws.Cells.Item(12,8).Value = 0.244689
$ws.Cells.Item(12,9).Value = 0.03578387091254728
$ws.Cells.Item(12,10).Value = 0.03578387091254728
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 1.683518666666667
$ws.Cells.Item(12,14).Value = 5.050556
$ws.Cells.Item(12,15).Value = 0.663513681459685
$ws.Cells.Item(12,16).Value = 0.663513681459685
$ws.Cells.Item(12,17).Value = 0.1373128330093333
$ws.Cells.Item(12,18).Value = 1.235815497084
$ws.Cells.Item(12,19).Value = 0.02374308792606238
$ws.Cells.Item(12,20).Value = 0.02374308792606238

$ws.Cells.Item(13,1).Value = "MuSCs"
$ws.Cells.Item(13,2).Value = "Efna5"
$ws.Cells.Item(13,3).Value = "Ephb6"
$ws.Cells.Item(13,4).Value = "Resolving-Mac"
$ws.Cells.Item(13,5).Value = 2
$ws.Cells.Item(13,6).Value = 0.6666666666666666
$ws.Cells.Item(13,7).Value = 0.081563
$ws.Cells.Item(13,8).Value = 0.244689
$ws.Cells.Item(13,9).Value = 0.03578387091254728
$ws.Cells.Item(13,10).Value = 0.03578387091254728
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 0.072253
$ws.Cells.Item(13,14).Value = 0.216759
$ws.Cells.Item(13,15).Value = 0.02847658002000569
$ws.Cells.Item(13,16).Value = 0.02847658002000569
$ws.Cells.Item(13,17).Value = 0.005893171439
$ws.Cells.Item(13,18).Value = 0.053038542951
$ws.Cells.Item(13,19).Value = 0.001019002263466707
$ws.Cells.Item(13,20).Value = 0.001019002263466707
